$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold font + thin border + centered) from AC1 (existing
# last header cell) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Team record values for all data rows (2 through 47)
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 70   # AD = column 30
    $ws.Cells.Item($row, 31).Value = 92   # AE = column 31
    $ws.Cells.Item($row, 32).Value = 0    # AF = column 32
}
